$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (date rolls from 05-06 to 05-07)
$ws.Name = "Through 2022-05-07"

# Update the "Total" column header label (shared string) for the current year
$ws.Range("I1").Value = "2022 (through 05-07)"

# Update the per-month Total (column I) values that changed with the new day's data
$ws.Range("I5").Value = 117
$ws.Range("I6").Value = 24
$ws.Range("I14").Value = 576
